$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-5 per regenerated save_data
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 7
